$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the worksheet (tab) from "Sheet1" to "MT_10"
$ws.Name = "MT_10"

# Update the active selection shown in the sheet view
$ws.Range("K60").Select()

# Each of the five 10-row blocks (rows 2-11, 13-22, 24-33, 35-44, 46-55) has a
# summary row (12, 23, 34, 45, 56) holding MIN(C) in column C and a literal
# max in column J. Column H originally referenced the workbook-wide $B$12 and
# column L referenced the workbook-wide $I$12; both now reference the
# corresponding per-block summary cells instead ($C$<end> and $J$<end>).
$blocks = @(
    @{Start = 2;  End = 12},
    @{Start = 13; End = 23},
    @{Start = 24; End = 34},
    @{Start = 35; End = 45},
    @{Start = 46; End = 56}
)

foreach ($blk in $blocks) {
    $start = $blk.Start
    $end = $blk.End
    $firstData = $start
    $lastData = $end - 1

    # Top row of the block (non-shared "master" formula)
    $ws.Range("H$firstData").Formula = "=(C$firstData-`$C`$$end)/100"
    $ws.Range("L$firstData").Formula = "=((F$firstData+G$firstData)/`$J`$$end)*100"

    # Remaining rows of the block, assigned together so the shared-formula
    # group is regenerated as a single block (matching the original layout).
    if ($lastData -gt $firstData) {
        $nextRow = $firstData + 1
        $ws.Range("H$nextRow`:H$lastData").Formula = "=(C$nextRow-`$C`$$end)/100"
        $ws.Range("L$nextRow`:L$lastData").Formula = "=((F$nextRow+G$nextRow)/`$J`$$end)*100"
    }
}
